$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark that currently sits in the
#    "Keepin it local Oreilly" paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Add the new run with the html5rocks URL to the last (empty)
#    paragraph of the document, just before the sectPr.
#    A trailing dummy character is used while inserting text/setting
#    language and while placing the new bookmark, because this COM
#    runtime mis-resolves a zero-length Range that sits exactly at
#    the last character position of a paragraph (i.e. right before
#    the paragraph mark). Appending (and later removing) one extra
#    character keeps every intermediate Range away from that
#    problematic boundary.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$rng = $lastPara.Range
$rng.InsertBefore("https://www.html5rocks.com/en/tutorials/offline/storage/X")
$rng.LanguageID = "en-US"

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bmPos = $lastPara.Range.End - 2

# ------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark right after the new text.
# ------------------------------------------------------------------
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 4. Remove the dummy placeholder character again.
# ------------------------------------------------------------------
$dummyRange = $d.Range($bmPos, $bmPos + 1)
$dummyRange.Delete()
